$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing record in row 4 (coordinates, address, phone) ---
$ws.Range("B4").Value = "'34.333483891624894"
$ws.Range("C4").Value = "'134.0444148832138"
$ws.Range("E4").Value = "'高松市中野町26-29"
$ws.Range("F4").Value = "'087-802-2067"

# --- Fill in previously-missing cells on row 47 ---
$ws.Range("A47").Value = "'46"
$ws.Range("H47").Value = "'"
$ws.Range("K47").Value = "'"
$ws.Range("L47").Value = "'"

# --- Append new row 48 ---
$ws.Range("A48").Value = "'47"
$ws.Range("B48").Value = "'34.34168814566888"
$ws.Range("C48").Value = "'134.04142199828937"
$ws.Range("D48").Value = "'地域密着型特別養護老人ホーム　マイルドハート番町"
$ws.Range("E48").Value = "'高松市番町三丁目14番15号"
$ws.Range("F48").Value = "'087-899-2710"
$ws.Range("G48").Value = "'https://www.utazufukushikai.or.jp/bancho/"
$ws.Range("I48").Value = "'地域密着型介護老人福祉施設入所者生活介護（地域密着型特別養護老人ホーム）"
$ws.Range("J48").Value = "'月火水木金土日"
$ws.Range("M48").Value = "'24時間対応"
$ws.Range("N48").Value = "'29"
